$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while preserving its existing
# number format (the cell displays a numeric-looking string but must be
# stored as text, not as a number).
function Set-TextValue($addr, $val, $fmt) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = $fmt
}

# Update the array-size values
$ws.Range("B3").Value = 50000
$ws.Range("B4").Value = 100000

# Update the Python runtime-time text values (Bubble Sort block)
Set-TextValue "D3" "94.06593" "[$-F400]h:mm:ss\ AM/PM"
Set-TextValue "D4" "381.88609" "[$]hh:mm:ss;@"
Set-TextValue "D5" "869.28647" "[$]hh:mm:ss;@"
Set-TextValue "D6" "1564.29378" "[$]hh:mm:ss;@"
Set-TextValue "D7" "2527.58961" "[$]hh:mm:ss;@"
Set-TextValue "D8" "3520.05186" "[$]hh:mm:ss;@"
Set-TextValue "D9" "4715.42362" "[$]hh:mm:ss;@"
Set-TextValue "D10" "7078.99730" "[$]hh:mm:ss;@"
Set-TextValue "D11" "9364.29803" "[$]hh:mm:ss;@"
Set-TextValue "D12" "11663.87637" "[$]hh:mm:ss;@"
Set-TextValue "D13" "13622.29164" "[$]hh:mm:ss;@"
Set-TextValue "D14" "16321.97137" "[$]hh:mm:ss;@"
Set-TextValue "D15" "19550.50262" "[$]hh:mm:ss;@"
Set-TextValue "D16" "24435.47438" "[$]hh:mm:ss;@"
Set-TextValue "D17" "25468.59240" "[$]hh:mm:ss;@"
Set-TextValue "D18" "29850.17472" "[$]hh:mm:ss;@"

# Update the selected cell to match the author's final cursor position
$null = $ws.Range("D10").Select()
